# Applies the "Found reliable way to select cells from table" commit:
#  - Tests!row17 gains a new data row (A17/B17/D17/E17) taken from the
#    "Test_Framework\Test_SearchRMNumber.xaml" / Success / Single result set.
#  - Tests sheet becomes the active (selected) tab, with the selection
#    parked just past the new data at E21.
#  - Result sheet is scrolled down (topLeftCell A4) with D8 selected, and
#    is no longer the active tab.

$wb = $excel.ActiveWorkbook

$tests  = $wb.Worksheets.Item("Tests")
$result = $wb.Worksheets.Item("Result")

# --- Tests: fill in the new row 17 (D17 stays blank, keeping its style) ---
$tests.Range("A17").Value = "Test_Framework\Test_SearchRMNumber.xaml"
$tests.Range("B17").Value = "Success"
$tests.Range("E17").Value = "Single result (Dec 2019)"

# --- Result: scroll so row 4 is at the top, and select D8 ---
$result.Activate() | Out-Null
$result.Range("D8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# --- Tests: make it the active/selected tab, selection resting at E21 ---
$tests.Activate() | Out-Null
$tests.Range("E21").Select() | Out-Null
